$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3
$ws.Cells.Item($row, 1).Value = 42632.880844907406
$ws.Cells.Item($row, 2).Value = 8
$ws.Cells.Item($row, 3).Value = "Buy"
$ws.Cells.Item($row, 4).Value = 12
$ws.Cells.Item($row, 5).Value = 14753
$ws.Cells.Item($row, 6).Value = 1499
$ws.Cells.Item($row, 7).Value = 58
$ws.Cells.Item($row, 8).Value = 37
$ws.Cells.Item($row, 9).Value = 74
$ws.Cells.Item($row, 10).Value = 24
$ws.Cells.Item($row, 11).Value = 8799
$ws.Cells.Item($row, 12).Value = 278
$ws.Cells.Item($row, 13).Value = 178
$ws.Cells.Item($row, 14).Value = 18
$ws.Cells.Item($row, 15).Value = 6
$ws.Cells.Item($row, 16).Value = "Bag"
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = 0.56000000000000005
$ws.Cells.Item($row, 19).Value = 0.10150000000000001
$ws.Cells.Item($row, 19).NumberFormat = $ws.Cells.Item(2, 19).NumberFormat
$ws.Cells.Item($row, 20).Value = -0.93
$ws.Cells.Item($row, 21).Value = 2.3199999999999998
$ws.Cells.Item($row, 22).Value = "N/A"
$ws.Cells.Item($row, 23).Value = 0
